$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.045207493600382
$ws.Cells.Item(2, 4).Value = 1.052815768352613
$ws.Cells.Item(2, 5).Value = 1.048872992161993
$ws.Cells.Item(2, 6).Value = 1.061448878614035
$ws.Cells.Item(2, 9).Value = 1.044975512436983
$ws.Cells.Item(2, 10).Value = 1.050268461179476
$ws.Cells.Item(2, 11).Value = 1.055563456390983
$ws.Cells.Item(2, 12).Value = 1.051631624664941
$ws.Cells.Item(2, 13).Value = 1.064172910069629
$ws.Cells.Item(2, 14).Value = 1.051759962861827
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.046153369513794
$ws.Cells.Item(3, 4).Value = 1.053573392722363
$ws.Cells.Item(3, 5).Value = 1.049767935389616
$ws.Cells.Item(3, 6).Value = 1.062337575659549
$ws.Cells.Item(3, 9).Value = 1.045230857899465
$ws.Cells.Item(3, 10).Value = 1.050861923263942
$ws.Cells.Item(3, 11).Value = 1.056134143469728
$ws.Cells.Item(3, 12).Value = 1.052338486476701
$ws.Cells.Item(3, 13).Value = 1.064876043494955
$ws.Cells.Item(3, 14).Value = 1.052354267730524
$ws.Cells.Item(4, 2).Value = 1.019999999999999
$ws.Cells.Item(4, 3).Value = 1.046765961332291
$ws.Cells.Item(4, 4).Value = 1.054064104069925
$ws.Cells.Item(4, 5).Value = 1.050347897833222
$ws.Cells.Item(4, 6).Value = 1.062913433190467
$ws.Cells.Item(4, 9).Value = 1.045395205484071
$ws.Cells.Item(4, 10).Value = 1.05124584450863
$ws.Cells.Item(4, 11).Value = 1.056503228756553
$ws.Cells.Item(4, 12).Value = 1.052796111583347
$ws.Cells.Item(4, 13).Value = 1.065331189175089
$ws.Cells.Item(4, 14).Value = 1.052738734187416
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.047023624815105
$ws.Cells.Item(5, 4).Value = 1.054270512109144
$ws.Cells.Item(5, 5).Value = 1.050591922003011
$ws.Cells.Item(5, 6).Value = 1.063155716429843
$ws.Cells.Item(5, 9).Value = 1.045464086448198
$ws.Cells.Item(5, 10).Value = 1.051407222859535
$ws.Cells.Item(5, 11).Value = 1.05665834631745
$ws.Cells.Item(5, 12).Value = 1.052988553089416
$ws.Cells.Item(5, 13).Value = 1.065522572232708
$ws.Cells.Item(5, 14).Value = 1.052900341714089
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.047066895234024
$ws.Cells.Item(6, 4).Value = 1.054305175512832
$ws.Cells.Item(6, 5).Value = 1.050632906890226
$ws.Cells.Item(6, 6).Value = 1.06319640811171
$ws.Cells.Item(6, 9).Value = 1.045475639489717
$ws.Cells.Item(6, 10).Value = 1.051434317668124
$ws.Cells.Item(6, 11).Value = 1.056684388517451
$ws.Cells.Item(6, 12).Value = 1.053020868111054
$ws.Cells.Item(6, 13).Value = 1.065554708589927
$ws.Cells.Item(6, 14).Value = 1.052927475000414
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.046769403735293
$ws.Cells.Item(7, 4).Value = 1.054066861662005
$ws.Cells.Item(7, 5).Value = 1.050351157680684
$ws.Cells.Item(7, 6).Value = 1.062916669835472
$ws.Cells.Item(7, 9).Value = 1.04539612670373
$ws.Cells.Item(7, 10).Value = 1.051248000942267
$ws.Cells.Item(7, 11).Value = 1.056505301626606
$ws.Cells.Item(7, 12).Value = 1.052798682777882
$ws.Cells.Item(7, 13).Value = 1.065333746289936
$ws.Cells.Item(7, 14).Value = 1.052740893683436
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.04552704294008
$ws.Cells.Item(8, 4).Value = 1.053071711311672
$ws.Cells.Item(8, 5).Value = 1.049175260562245
$ws.Cells.Item(8, 6).Value = 1.061749049122144
$ws.Cells.Item(8, 9).Value = 1.045061989087975
$ws.Cells.Item(8, 10).Value = 1.050469042187814
$ws.Cells.Item(8, 11).Value = 1.055756360948799
$ws.Cells.Item(8, 12).Value = 1.051870461915411
$ws.Cells.Item(8, 13).Value = 1.064410501101488
$ws.Cells.Item(8, 14).Value = 1.051960828718203
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.043342078546998
$ws.Cells.Item(9, 4).Value = 1.051321848907006
$ws.Cells.Item(9, 5).Value = 1.047109935987333
$ws.Cells.Item(9, 6).Value = 1.059697828512584
$ws.Cells.Item(9, 9).Value = 1.044466496455923
$ws.Cells.Item(9, 10).Value = 1.049095783357623
$ws.Cells.Item(9, 11).Value = 1.054435241792597
$ws.Cells.Item(9, 12).Value = 1.050236692620307
$ws.Cells.Item(9, 13).Value = 1.062784990609511
$ws.Cells.Item(9, 14).Value = 1.050585619702959
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.041888343035839
$ws.Cells.Item(10, 4).Value = 1.050157856477263
$ws.Cells.Item(10, 5).Value = 1.045737675842625
$ws.Cells.Item(10, 6).Value = 1.058334653721522
$ws.Cells.Item(10, 9).Value = 1.04406503044227
$ws.Cells.Item(10, 10).Value = 1.048179908778896
$ws.Cells.Item(10, 11).Value = 1.053553621809261
$ws.Cells.Item(10, 12).Value = 1.049148839854411
$ws.Cells.Item(10, 13).Value = 1.061702306274827
$ws.Cells.Item(10, 14).Value = 1.049668444477278
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.041259561541669
$ws.Cells.Item(11, 4).Value = 1.049654464496067
$ws.Cells.Item(11, 5).Value = 1.045144585259859
$ws.Cells.Item(11, 6).Value = 1.057745423127409
$ws.Cells.Item(11, 9).Value = 1.04389013771616
$ws.Cells.Item(11, 10).Value = 1.047783249471658
$ws.Cells.Item(11, 11).Value = 1.053171675974876
$ws.Cells.Item(11, 12).Value = 1.048678116104749
$ws.Cells.Item(11, 13).Value = 1.061233741501652
$ws.Cells.Item(11, 14).Value = 1.049271221868326
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.041026109555165
$ws.Cells.Item(12, 4).Value = 1.049467577469318
$ws.Cells.Item(12, 5).Value = 1.04492445266633
$ws.Cells.Item(12, 6).Value = 1.05752671328158
$ws.Cells.Item(12, 9).Value = 1.043825016756297
$ws.Cells.Item(12, 10).Value = 1.047635901522475
$ws.Cells.Item(12, 11).Value = 1.053029775443789
$ws.Cells.Item(12, 12).Value = 1.048503318059982
$ws.Cells.Item(12, 13).Value = 1.06105973361404
$ws.Cells.Item(12, 14).Value = 1.049123664668154
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.041076181011207
$ws.Cells.Item(13, 4).Value = 1.049507661042597
$ws.Cells.Item(13, 5).Value = 1.044971664246809
$ws.Cells.Item(13, 6).Value = 1.057573620175134
$ws.Cells.Item(13, 9).Value = 1.043838992565162
$ws.Cells.Item(13, 10).Value = 1.047667508641482
$ws.Cells.Item(13, 11).Value = 1.053060214878649
$ws.Cells.Item(13, 12).Value = 1.048540810566788
$ws.Cells.Item(13, 13).Value = 1.061097057161301
$ws.Cells.Item(13, 14).Value = 1.049155316672895
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.041240262170806
$ws.Cells.Item(14, 4).Value = 1.049639014406879
$ws.Cells.Item(14, 5).Value = 1.045126385603199
$ws.Cells.Item(14, 6).Value = 1.057727341299839
$ws.Cells.Item(14, 9).Value = 1.043884758017661
$ws.Cells.Item(14, 10).Value = 1.047771069863635
$ws.Cells.Item(14, 11).Value = 1.053159947015346
$ws.Cells.Item(14, 12).Value = 1.048663666223382
$ws.Cells.Item(14, 13).Value = 1.061219357174802
$ws.Cells.Item(14, 14).Value = 1.049259024963862
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.041341371991329
$ws.Cells.Item(15, 4).Value = 1.049719958201168
$ws.Cells.Item(15, 5).Value = 1.045221736793073
$ws.Cells.Item(15, 6).Value = 1.057822074756106
$ws.Cells.Item(15, 9).Value = 1.043912934700924
$ws.Cells.Item(15, 10).Value = 1.047834875918104
$ws.Cells.Item(15, 11).Value = 1.053221391484141
$ws.Cells.Item(15, 12).Value = 1.04873936826839
$ws.Cells.Item(15, 13).Value = 1.061294715310936
$ws.Cells.Item(15, 14).Value = 1.049322921630248
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.041930087929078
$ws.Cells.Item(16, 4).Value = 1.050191278216775
$ws.Cells.Item(16, 5).Value = 1.045777060749973
$ws.Cells.Item(16, 6).Value = 1.058373780902807
$ws.Cells.Item(16, 9).Value = 1.044076615308
$ws.Cells.Item(16, 10).Value = 1.048206232153899
$ws.Cells.Item(16, 11).Value = 1.053578966197314
$ws.Cells.Item(16, 12).Value = 1.049180087151526
$ws.Cells.Item(16, 13).Value = 1.061733408630788
$ws.Cells.Item(16, 14).Value = 1.049694805234493
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.042299560780446
$ws.Cells.Item(17, 4).Value = 1.050487092943141
$ws.Cells.Item(17, 5).Value = 1.046125698137328
$ws.Cells.Item(17, 6).Value = 1.058720129183019
$ws.Cells.Item(17, 9).Value = 1.04417900567087
$ws.Cells.Item(17, 10).Value = 1.048439153464285
$ws.Cells.Item(17, 11).Value = 1.053803210904849
$ws.Cells.Item(17, 12).Value = 1.049456625886643
$ws.Cells.Item(17, 13).Value = 1.062008655695712
$ws.Cells.Item(17, 14).Value = 1.049928057319856
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.042515135092259
$ws.Cells.Item(18, 4).Value = 1.05065969682492
$ws.Cells.Item(18, 5).Value = 1.046329159158517
$ws.Cells.Item(18, 6).Value = 1.058922247857415
$ws.Cells.Item(18, 9).Value = 1.044238626380951
$ws.Cells.Item(18, 10).Value = 1.048575004783753
$ws.Cells.Item(18, 11).Value = 1.0539339898238
$ws.Cells.Item(18, 12).Value = 1.04961795740931
$ws.Cells.Item(18, 13).Value = 1.062169226182141
$ws.Cells.Item(18, 14).Value = 1.050064101563781
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.042588651705299
$ws.Cells.Item(19, 4).Value = 1.05071856045049
$ws.Cells.Item(19, 5).Value = 1.046398552172107
$ws.Cells.Item(19, 6).Value = 1.058991181949795
$ws.Cells.Item(19, 9).Value = 1.044258938209633
$ws.Cells.Item(19, 10).Value = 1.048621325231755
$ws.Cells.Item(19, 11).Value = 1.05397857876455
$ws.Cells.Item(19, 12).Value = 1.049672972541146
$ws.Cells.Item(19, 13).Value = 1.062223980524072
$ws.Cells.Item(19, 14).Value = 1.050110487792132
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.042259912898079
$ws.Cells.Item(20, 4).Value = 1.050455348590473
$ws.Cells.Item(20, 5).Value = 1.046088281595802
$ws.Cells.Item(20, 6).Value = 1.058682958979809
$ws.Cells.Item(20, 9).Value = 1.044168030685577
$ws.Cells.Item(20, 10).Value = 1.048414164005873
$ws.Cells.Item(20, 11).Value = 1.053779153551337
$ws.Cells.Item(20, 12).Value = 1.049426952666915
$ws.Cells.Item(20, 13).Value = 1.061979121848248
$ws.Cells.Item(20, 14).Value = 1.049903032373547
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.041191941440096
$ws.Cells.Item(21, 4).Value = 1.049600331476979
$ws.Cells.Item(21, 5).Value = 1.045080819381287
$ws.Cells.Item(21, 6).Value = 1.057682069924756
$ws.Cells.Item(21, 9).Value = 1.043871285590259
$ws.Cells.Item(21, 10).Value = 1.047740573961637
$ws.Cells.Item(21, 11).Value = 1.053130579175068
$ws.Cells.Item(21, 12).Value = 1.048627486913505
$ws.Cells.Item(21, 13).Value = 1.061183341812555
$ws.Cells.Item(21, 14).Value = 1.049228485754186
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.04052107601516
$ws.Cells.Item(22, 4).Value = 1.049063299570379
$ws.Cells.Item(22, 5).Value = 1.044448358566002
$ws.Cells.Item(22, 6).Value = 1.057053678182802
$ws.Cells.Item(22, 9).Value = 1.043683796214431
$ws.Cells.Item(22, 10).Value = 1.047316997864724
$ws.Cells.Item(22, 11).Value = 1.05272262834333
$ws.Cells.Item(22, 12).Value = 1.048125119237862
$ws.Cells.Item(22, 13).Value = 1.060683223363567
$ws.Cells.Item(22, 14).Value = 1.048804308130631
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.040876656270166
$ws.Cells.Item(23, 4).Value = 1.049347937636192
$ws.Cells.Item(23, 5).Value = 1.044783545559926
$ws.Cells.Item(23, 6).Value = 1.057386714092445
$ws.Cells.Item(23, 9).Value = 1.043783274334225
$ws.Cells.Item(23, 10).Value = 1.047541549255902
$ws.Cells.Item(23, 11).Value = 1.052938906237631
$ws.Cells.Item(23, 12).Value = 1.048391406168908
$ws.Cells.Item(23, 13).Value = 1.060948324366908
$ws.Cells.Item(23, 14).Value = 1.049029178410541
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.042277827862392
$ws.Cells.Item(24, 4).Value = 1.050469692311295
$ws.Cells.Item(24, 5).Value = 1.046105188190692
$ws.Cells.Item(24, 6).Value = 1.058699754287123
$ws.Cells.Item(24, 9).Value = 1.044172990123784
$ws.Cells.Item(24, 10).Value = 1.048425455690471
$ws.Cells.Item(24, 11).Value = 1.05379002409335
$ws.Cells.Item(24, 12).Value = 1.049440360621474
$ws.Cells.Item(24, 13).Value = 1.061992466850059
$ws.Cells.Item(24, 14).Value = 1.049914340093632
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.043906436227076
$ws.Cells.Item(25, 4).Value = 1.051773781076882
$ws.Cells.Item(25, 5).Value = 1.047643062963981
$ws.Cells.Item(25, 6).Value = 1.060227365603041
$ws.Cells.Item(25, 9).Value = 1.044621236063268
$ws.Cells.Item(25, 10).Value = 1.049450872889677
$ws.Cells.Item(25, 11).Value = 1.054776941238725
$ws.Cells.Item(25, 12).Value = 1.050658832169495
$ws.Cells.Item(25, 13).Value = 1.063205054883897
$ws.Cells.Item(25, 14).Value = 1.050941213502878
